$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.744.97'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '3.315.95'
$ws.Range('E3').Value = '  +2.40%  '
$ws.Range('D5').Value = '604.82'
$ws.Range('E5').Value = '  +1.64%  '
$ws.Range('D6').Value = '141.70'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D8').Value = '3.315.21'
$ws.Range('E8').Value = '  +2.53%  '
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('E10').Value = '  +1.48%  '
$ws.Range('D11').Value = '5.54'
$ws.Range('E11').Value = '  +3.72%  '
$ws.Range('D12').Value = '0.469'
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('D14').Value = '34.98'
$ws.Range('E14').Value = '  +1.44%  '
$ws.Range('D15').Value = '3.863.56'
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range('D16').Value = '0.121'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').Value = '3.318.58'
$ws.Range('E17').Value = '  +2.45%  '
$ws.Range('D18').Value = '63.833.49'
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('D19').Value = '6.86'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('D20').Value = '481.48'
$ws.Range('E20').Value = '  +1.27%  '
$ws.Range('D21').Value = '14.09'
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').Value = '0.739'
$ws.Range('E22').Value = '  +1.54%  '
$ws.Range('D23').Value = '8.00'
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('D24').Value = '14.07'
$ws.Range('E24').Value = '  +6.66%  '
$ws.Range('D25').Value = '85.14'
$ws.Range('E25').Value = '  +1.35%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').Value = '2.77'
$ws.Range('E27').Value = '  +1.29%  '
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('D29').Value = '8.24'
$ws.Range('E29').Value = '  +1.60%  '
$ws.Range('D30').Value = '7.18'
$ws.Range('E30').Value = '  -5.39%  '
$ws.Range('D31').Value = '2.15'
$ws.Range('E31').Value = '  +2.10%  '
$ws.Range('E32').Value = '  +5.42%  '
$ws.Range('D33').Value = '0.106'
$ws.Range('E33').Value = '  -1.38%  '
$ws.Range('D34').Value = '2.52'
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('D35').Value = '1.10'
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('D36').Value = '6.08'
$ws.Range('E36').Value = '  +2.85%  '
$ws.Range('D37').Value = '52.44'
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('E38').Value = '  +4.82%  '
$ws.Range('E39').Value = '  +1.53%  '
$ws.Range('D40').Value = '433.89'
$ws.Range('E40').Value = '  +2.59%  '
$ws.Range('D41').Value = '3.123.29'
$ws.Range('E41').Value = '  +4.95%  '
$ws.Range('D42').Value = '0.119'
$ws.Range('E42').Value = '  +7.74%  '
$ws.Range('D43').Value = '8.36'
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range('D44').Value = '2.74'
$ws.Range('E44').Value = '  -0.98%  '
$ws.Range('D45').Value = '0.267'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').Value = '2.24'
$ws.Range('E46').Value = '  +3.37%  '
$ws.Range('D47').Value = '36.92'
$ws.Range('E47').Value = '  +7.54%  '
$ws.Range('D48').Value = '26.38'
$ws.Range('E48').Value = '  +1.58%  '
$ws.Range('E50').Value = '  -3.13%  '
$ws.Range('E51').Value = '  -0.67%  '
